$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2 and 3 with the new capital structure figures
foreach ($r in 2,3) {
    $ws.Range("K$r").Value = -0.036

    $ws.Range("W$r").Value = 0.175609756097561
    $ws.Range("X$r").Value = 0.1022842077691
    $ws.Range("Y$r").Value = 0.07332554832846101

    $ws.Range("AA$r").Value = 1.25925925925926
    $ws.Range("AB$r").Value = 0.08687786186139404
    $ws.Range("AC$r").Value = 1.172381397397865

    $ws.Range("AD$r").Value = 0.217
    $ws.Range("AF$r").Value = 0.217
    $ws.Range("AG$r").Value = 0.217

    $ws.Range("AH$r").Value = 0.2916666666666667
    $ws.Range("AI$r").Value = -9.04166666666667
    $ws.Range("AJ$r").Value = 0.2916666666666667
    $ws.Range("AK$r").Value = -9.04166666666667
}
